$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three rows 25, 26, 27 are cyclically rotated:
# new row 25 = old row 26 data
# new row 26 = old row 27 data
# new row 27 = old row 25 data
# (column D "Rodlistade" stays "NT" for all three, so it is unaffected)

# Capture old values before overwriting anything (use Value2, since Value is unreliable here)
$old25 = @{
    A = $ws.Range("A25").Value2
    B = $ws.Range("B25").Value2
    E = $ws.Range("E25").Value2
    F = $ws.Range("F25").Value2
    G = $ws.Range("G25").Value2
    H = $ws.Range("H25").Value2
    Q = $ws.Range("Q25").Value2
    R = $ws.Range("R25").Value2
    S = $ws.Range("S25").Value2
    Z = $ws.Range("Z25").Value2
    AB = $ws.Range("AB25").Value2
}

$old26 = @{
    A = $ws.Range("A26").Value2
    B = $ws.Range("B26").Value2
    E = $ws.Range("E26").Value2
    F = $ws.Range("F26").Value2
    G = $ws.Range("G26").Value2
    H = $ws.Range("H26").Value2
    Q = $ws.Range("Q26").Value2
    R = $ws.Range("R26").Value2
    S = $ws.Range("S26").Value2
    Z = $ws.Range("Z26").Value2
    AB = $ws.Range("AB26").Value2
}

$old27 = @{
    A = $ws.Range("A27").Value2
    B = $ws.Range("B27").Value2
    E = $ws.Range("E27").Value2
    F = $ws.Range("F27").Value2
    G = $ws.Range("G27").Value2
    H = $ws.Range("H27").Value2
    Q = $ws.Range("Q27").Value2
    R = $ws.Range("R27").Value2
    S = $ws.Range("S27").Value2
    Z = $ws.Range("Z27").Value2
    AB = $ws.Range("AB27").Value2
}

# Write old26 values into row 25
$ws.Range("A25").Value2 = $old26.A
$ws.Range("B25").Value2 = $old26.B
$ws.Range("E25").Value2 = $old26.E
$ws.Range("F25").Value2 = $old26.F
$ws.Range("G25").Value2 = $old26.G
$ws.Range("H25").Value2 = $old26.H
$ws.Range("Q25").Value2 = $old26.Q
$ws.Range("R25").Value2 = $old26.R
$ws.Range("S25").Value2 = $old26.S
$ws.Range("Z25").Value2 = $old26.Z
$ws.Range("AB25").Value2 = $old26.AB

# Write old27 values into row 26
$ws.Range("A26").Value2 = $old27.A
$ws.Range("B26").Value2 = $old27.B
$ws.Range("E26").Value2 = $old27.E
$ws.Range("F26").Value2 = $old27.F
$ws.Range("G26").Value2 = $old27.G
$ws.Range("H26").Value2 = $old27.H
$ws.Range("Q26").Value2 = $old27.Q
$ws.Range("R26").Value2 = $old27.R
$ws.Range("S26").Value2 = $old27.S
$ws.Range("Z26").Value2 = $old27.Z
$ws.Range("AB26").Value2 = $old27.AB

# Write old25 values into row 27
$ws.Range("A27").Value2 = $old25.A
$ws.Range("B27").Value2 = $old25.B
$ws.Range("E27").Value2 = $old25.E
$ws.Range("F27").Value2 = $old25.F
$ws.Range("G27").Value2 = $old25.G
$ws.Range("H27").Value2 = $old25.H
$ws.Range("Q27").Value2 = $old25.Q
$ws.Range("R27").Value2 = $old25.R
$ws.Range("S27").Value2 = $old25.S
$ws.Range("Z27").Value2 = $old25.Z
$ws.Range("AB27").Value2 = $old25.AB
